# Auto-generated script to apply scheduled market-data refresh to Moogle_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) per sheet with refreshed values
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 594.5263
$ws.Range("I9").Value = 318.46155
$ws.Range("J9").Value = 1192.6666
$ws.Range("K9").Value = 318.46155
$ws.Range("L9").Value = 1192.6666
$ws.Range("M9").Value = -149.46155
$ws.Range("N9").Value = -1530.6666
$ws.Range("H28").Value = 1982.5
$ws.Range("I28").Value = 1896
$ws.Range("J28").Value = 2198.75
$ws.Range("K28").Value = 1896
$ws.Range("L28").Value = 2198.75
$ws.Range("M28").Value = -1411
$ws.Range("N28").Value = -3168.75
$ws.Range("H76").Value = 6736.1816
$ws.Range("J76").Value = 7514.2856
$ws.Range("L76").Value = 7514.2856
$ws.Range("N76").Value = -8144.2856
$ws.Range("H79").Value = 6736.1816
$ws.Range("J79").Value = 7514.2856
$ws.Range("L79").Value = 7514.2856
$ws.Range("N79").Value = -9698.285599999999
$ws.Range("H86").Value = 7175.7744
$ws.Range("I86").Value = 6688.619
$ws.Range("J86").Value = 8198.799999999999
$ws.Range("K86").Value = 6688.619
$ws.Range("L86").Value = 8198.799999999999
$ws.Range("M86").Value = -5565.619
$ws.Range("N86").Value = -10444.8
$ws.Range("H89").Value = 7175.7744
$ws.Range("I89").Value = 6688.619
$ws.Range("J89").Value = 8198.799999999999
$ws.Range("K89").Value = 33443.095
$ws.Range("L89").Value = 40994
$ws.Range("M89").Value = -27827.095
$ws.Range("N89").Value = -52226
$ws.Range("H101").Value = 906.625
$ws.Range("I101").Value = 893.3333
$ws.Range("J101").Value = 946.5
$ws.Range("K101").Value = 2679.9999
$ws.Range("L101").Value = 2839.5
$ws.Range("M101").Value = -1057.9999
$ws.Range("N101").Value = -6083.5
$ws.Range("H113").Value = 5225.68
$ws.Range("I113").Value = 4191.316
$ws.Range("K113").Value = 4191.316
$ws.Range("M113").Value = -937.3159999999998
$ws.Range("H125").Value = 72846.21000000001
$ws.Range("I125").Value = 2349.75
$ws.Range("J125").Value = 101044.8
$ws.Range("K125").Value = 21147.75
$ws.Range("L125").Value = 909403.2000000001
$ws.Range("M125").Value = -18687.75
$ws.Range("N125").Value = -914323.2000000001
$ws.Range("H137").Value = 2418.516
$ws.Range("I137").Value = 2099.074
$ws.Range("J137").Value = 4574.75
$ws.Range("K137").Value = 6297.222
$ws.Range("L137").Value = 13724.25
$ws.Range("M137").Value = -3747.222
$ws.Range("N137").Value = -18824.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9777.361000000001
$ws.Range("I32").Value = 3777.3225
$ws.Range("J32").Value = 46977.6
$ws.Range("K32").Value = 3777.3225
$ws.Range("L32").Value = 46977.6
$ws.Range("M32").Value = -3490.3225
$ws.Range("N32").Value = -47551.6
$ws.Range("H74").Value = 5257.5
$ws.Range("I74").Value = 1456.8462
$ws.Range("K74").Value = 1456.8462
$ws.Range("M74").Value = -582.8462
$ws.Range("H77").Value = 5257.5
$ws.Range("I77").Value = 1456.8462
$ws.Range("K77").Value = 7284.231
$ws.Range("M77").Value = -2916.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10246.714
$ws.Range("I31").Value = 3977.2856
$ws.Range("J31").Value = 22785.572
$ws.Range("K31").Value = 3977.2856
$ws.Range("L31").Value = 22785.572
$ws.Range("M31").Value = -3682.2856
$ws.Range("N31").Value = -23375.572
$ws.Range("H34").Value = 10246.714
$ws.Range("I34").Value = 3977.2856
$ws.Range("J34").Value = 22785.572
$ws.Range("K34").Value = 3977.2856
$ws.Range("L34").Value = 22785.572
$ws.Range("M34").Value = -3775.2856
$ws.Range("N34").Value = -23189.572
$ws.Range("H41").Value = 16280.75
$ws.Range("I41").Value = 7483.75
$ws.Range("K41").Value = 7483.75
$ws.Range("M41").Value = -7055.75
$ws.Range("H44").Value = 5500
$ws.Range("I44").Value = 6000
$ws.Range("J44").Value = 5000
$ws.Range("K44").Value = 6000
$ws.Range("L44").Value = 5000
$ws.Range("M44").Value = -5558
$ws.Range("N44").Value = -5884
$ws.Range("H50").Value = 49573.125
$ws.Range("J50").Value = 68797.39999999999
$ws.Range("L50").Value = 68797.39999999999
$ws.Range("N50").Value = -70047.39999999999
$ws.Range("H60").Value = 11498.823
$ws.Range("I60").Value = 4535.0713
$ws.Range("J60").Value = 43996.332
$ws.Range("K60").Value = 4535.0713
$ws.Range("L60").Value = 43996.332
$ws.Range("M60").Value = -4024.0713
$ws.Range("N60").Value = -45018.332
$ws.Range("H107").Value = 1490.8334
$ws.Range("I107").Value = 1229.4
$ws.Range("K107").Value = 1229.4
$ws.Range("M107").Value = 690.5999999999999
$ws.Range("H134").Value = 6263.5293
$ws.Range("I134").Value = 5105.9287
$ws.Range("J134").Value = 11665.667
$ws.Range("K134").Value = 15317.7861
$ws.Range("L134").Value = 34997.001
$ws.Range("M134").Value = -12782.7861
$ws.Range("N134").Value = -40067.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 139989
$ws.Range("J37").Value = 139989
$ws.Range("L37").Value = 419967
$ws.Range("N37").Value = -420191
$ws.Range("H113").Value = 630.4375
$ws.Range("I113").Value = 397.66666
$ws.Range("J113").Value = 684.1539
$ws.Range("K113").Value = 1192.99998
$ws.Range("L113").Value = 2052.4617
$ws.Range("M113").Value = 977.0000199999999
$ws.Range("N113").Value = -6392.4617

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1727
$ws.Range("I113").Value = 1638.2858
$ws.Range("K113").Value = 1638.2858
$ws.Range("M113").Value = 531.7141999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9089
$ws.Range("I132").Value = 3787.5
$ws.Range("K132").Value = 11362.5
$ws.Range("M132").Value = -8832.5
$ws.Range("H136").Value = 4691.1396
$ws.Range("I136").Value = 1793.2963
$ws.Range("K136").Value = 5379.8889
$ws.Range("M136").Value = -2829.8889

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 17060
$ws.Range("J74").Value = 17060
$ws.Range("L74").Value = 17060
$ws.Range("N74").Value = -18932
$ws.Range("H77").Value = 17060
$ws.Range("J77").Value = 17060
$ws.Range("L77").Value = 51180
$ws.Range("N77").Value = -60540
$ws.Range("H100").Value = 804
$ws.Range("J100").Value = 915.6667
$ws.Range("L100").Value = 1831.3334
$ws.Range("N100").Value = -2913.3334
$ws.Range("H122").Value = 3471.5386
$ws.Range("I122").Value = 3398.6
$ws.Range("J122").Value = 3714.6667
$ws.Range("K122").Value = 10195.8
$ws.Range("L122").Value = 11144.0001
$ws.Range("M122").Value = -7745.799999999999
$ws.Range("N122").Value = -16044.0001
$ws.Range("H132").Value = 2924.7585
$ws.Range("I132").Value = 2308.074
$ws.Range("J132").Value = 11250
$ws.Range("K132").Value = 6924.222
$ws.Range("L132").Value = 33750
$ws.Range("M132").Value = -4394.222
$ws.Range("N132").Value = -38810
